$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038396418900487
$ws.Range("D2").Value = 1.039620685560593
$ws.Range("E2").Value = 1.037056314358973
$ws.Range("I2").Value = 1.037963303409046
$ws.Range("J2").Value = 1.043494083145059
$ws.Range("K2").Value = 1.042405354590797
$ws.Range("L2").Value = 1.039848288949243
$ws.Range("N2").Value = 1.044975964433569
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03941153802151
$ws.Range("D3").Value = 1.040363030051888
$ws.Range("E3").Value = 1.037921293360791
$ws.Range("I3").Value = 1.038210639212956
$ws.Range("J3").Value = 1.044153648859432
$ws.Range("K3").Value = 1.042958134089813
$ws.Range("L3").Value = 1.040522850538266
$ws.Range("N3").Value = 1.04563646680691
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.040068507712754
$ws.Range("D4").Value = 1.040843442192237
$ws.Range("E4").Value = 1.038481451812584
$ws.Range("I4").Value = 1.038369553975853
$ws.Range("J4").Value = 1.044579961597173
$ws.Range("K4").Value = 1.0433152197141
$ws.Range("L4").Value = 1.040959157328005
$ws.Range("N4").Value = 1.046063384957637
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.040344726103056
$ws.Range("D5").Value = 1.041045422223164
$ws.Range("E5").Value = 1.038717051488431
$ws.Range("I5").Value = 1.038436091343087
$ws.Range("J5").Value = 1.044759070276224
$ws.Range("K5").Value = 1.04346519432889
$ws.Range("L5").Value = 1.041142537036599
$ws.Range("N5").Value = 1.046242747991557
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.040391106037989
$ws.Range("D6").Value = 1.041079336397678
$ws.Range("E6").Value = 1.03875661608089
$ws.Range("I6").Value = 1.038447247393526
$ws.Range("J6").Value = 1.044789136756702
$ws.Range("K6").Value = 1.043490367253785
$ws.Range("L6").Value = 1.041173324717005
$ws.Range("N6").Value = 1.046272857169885
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.040072198445487
$ws.Range("D7").Value = 1.040846141002456
$ws.Range("E7").Value = 1.038484599479533
$ws.Range("I7").Value = 1.038370444114322
$ws.Range("J7").Value = 1.044582355301624
$ws.Range("K7").Value = 1.043317224249554
$ws.Range("L7").Value = 1.040961607829042
$ws.Range("N7").Value = 1.046065782061423
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038739458165562
$ws.Range("D8").Value = 1.039871550205605
$ws.Range("E8").Value = 1.037348542459514
$ws.Range("I8").Value = 1.038047125063452
$ws.Range("J8").Value = 1.043717083510509
$ws.Range("K8").Value = 1.042592292805682
$ws.Range("L8").Value = 1.040076297004542
$ws.Range("N8").Value = 1.045199281485117
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036391931012921
$ws.Range("D9").Value = 1.038154740561936
$ws.Range("E9").Value = 1.035350219933605
$ws.Range("I9").Value = 1.03746877123346
$ws.Range("J9").Value = 1.042188787090197
$ws.Range("K9").Value = 1.041310305283455
$ws.Range("L9").Value = 1.038514918835396
$ws.Range("N9").Value = 1.04366881470859
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.034827556100441
$ws.Range("D10").Value = 1.037010620844664
$ws.Range("E10").Value = 1.034020445078015
$ws.Range("I10").Value = 1.03707742307684
$ws.Range("J10").Value = 1.041167550810085
$ws.Range("K10").Value = 1.040452612471672
$ws.Range("L10").Value = 1.037473131310368
$ws.Range("N10").Value = 1.042646128155822
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034150319775399
$ws.Range("D11").Value = 1.036515315044365
$ws.Range("E11").Value = 1.033445226682285
$ws.Range("I11").Value = 1.036906598186426
$ws.Range("J11").Value = 1.04072478805316
$ws.Range("K11").Value = 1.040080509291741
$ws.Range("L11").Value = 1.037021827118311
$ws.Range("N11").Value = 1.042202736624997
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.033898786451952
$ws.Range("D12").Value = 1.036331353118047
$ws.Range("E12").Value = 1.033231653281958
$ws.Range("I12").Value = 1.036842940872147
$ws.Range("J12").Value = 1.040560242460601
$ws.Range("K12").Value = 1.039942186458222
$ws.Range("L12").Value = 1.036854162633594
$ws.Range("N12").Value = 1.042037957358823
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.033952740188065
$ws.Range("D13").Value = 1.036370812806503
$ws.Range("E13").Value = 1.033277461492454
$ws.Range("I13").Value = 1.036856604870111
$ws.Range("J13").Value = 1.040595541833245
$ws.Range("K13").Value = 1.039971862015948
$ws.Range("L13").Value = 1.036890128584332
$ws.Range("N13").Value = 1.042073306860626
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034129527494746
$ws.Range("D14").Value = 1.036500108348603
$ws.Range("E14").Value = 1.033427570832279
$ws.Range("I14").Value = 1.03690134044038
$ws.Range("J14").Value = 1.040711188369849
$ws.Range("K14").Value = 1.040069077665961
$ws.Range("L14").Value = 1.037007968528883
$ws.Range("N14").Value = 1.042189117628576
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.034238454961074
$ws.Range("D15").Value = 1.036579773837525
$ws.Range("E15").Value = 1.033520069876935
$ws.Range("I15").Value = 1.036928876299597
$ws.Range("J15").Value = 1.040782430928257
$ws.Range("K15").Value = 1.040128961253905
$ws.Range("L15").Value = 1.0370805696468
$ws.Range("N15").Value = 1.042260461359589
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.034872505179979
$ws.Range("D16").Value = 1.037043494955556
$ws.Range("E16").Value = 1.034058632773442
$ws.Range("I16").Value = 1.037088731350992
$ws.Range("J16").Value = 1.041196923740924
$ws.Range("K16").Value = 1.040477292680928
$ws.Range("L16").Value = 1.037503078645048
$ws.Range("N16").Value = 1.042675542799592
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035270267925458
$ws.Range("D17").Value = 1.037334403822521
$ws.Range("E17").Value = 1.034396615636612
$ws.Range("I17").Value = 1.037188637987529
$ws.Range("J17").Value = 1.041456774331181
$ws.Range("K17").Value = 1.040695600303694
$ws.Range("L17").Value = 1.037768053319139
$ws.Range("N17").Value = 1.042935762407493
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.035502290509566
$ws.Range("D18").Value = 1.037504096174376
$ws.Range("E18").Value = 1.034593811512819
$ws.Range("I18").Value = 1.037246779774598
$ws.Range("J18").Value = 1.041608286653747
$ws.Range("K18").Value = 1.040822866277634
$ws.Range("L18").Value = 1.03792258891046
$ws.Range("N18").Value = 1.043087489894935
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035581406637516
$ws.Range("D19").Value = 1.037561958574075
$ws.Range("E19").Value = 1.034661059718616
$ws.Range("I19").Value = 1.037266582213094
$ws.Range("J19").Value = 1.041659939245188
$ws.Range("K19").Value = 1.040866248938241
$ws.Range("L19").Value = 1.037975278203827
$ws.Range("N19").Value = 1.043139215838981
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035227590266761
$ws.Range("D20").Value = 1.037303191000171
$ws.Range("E20").Value = 1.034360347456521
$ws.Range("I20").Value = 1.037177932608958
$ws.Range("J20").Value = 1.041428900417702
$ws.Range("K20").Value = 1.040672185112575
$ws.Range("L20").Value = 1.037739626073207
$ws.Range("N20").Value = 1.04290784890986
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034077467419894
$ws.Range("D21").Value = 1.036462033566298
$ws.Range("E21").Value = 1.033383364928893
$ws.Range("I21").Value = 1.036888172596296
$ws.Range("J21").Value = 1.040677135656652
$ws.Range("K21").Value = 1.040040453033717
$ws.Range("L21").Value = 1.036973268426535
$ws.Range("N21").Value = 1.04215501655662
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033354469206668
$ws.Range("D22").Value = 1.035933261545139
$ws.Range("E22").Value = 1.032769608620865
$ws.Range("I22").Value = 1.036704801051549
$ws.Range("J22").Value = 1.040203986750297
$ws.Range("K22").Value = 1.039642638347543
$ws.Range("L22").Value = 1.036491255572686
$ws.Range("N22").Value = 1.041681195724547
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.033737731908398
$ws.Range("D23").Value = 1.036213564120275
$ws.Range("E23").Value = 1.033094923728712
$ws.Range("I23").Value = 1.036802122289192
$ws.Range("J23").Value = 1.040454857718603
$ws.Range("K23").Value = 1.039853586006622
$ws.Range("L23").Value = 1.036746795963414
$ws.Range("N23").Value = 1.041932422958403
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035246874419378
$ws.Range("D24").Value = 1.037317294699787
$ws.Range("E24").Value = 1.03437673531236
$ws.Range("I24").Value = 1.037182770316772
$ws.Range("J24").Value = 1.041441495606614
$ws.Range("K24").Value = 1.040682765643264
$ws.Range("L24").Value = 1.037752471182726
$ws.Range("N24").Value = 1.042920461985385
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.036998710931856
$ws.Range("D25").Value = 1.038598506945254
$ws.Range("E25").Value = 1.035866407762013
$ws.Range("I25").Value = 1.037619309577154
$ws.Range("J25").Value = 1.042584308842241
$ws.Range("K25").Value = 1.04164226721168
$ws.Range("L25").Value = 1.038918728359526
$ws.Range("N25").Value = 1.044064898146888
